$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing shared string value (row 1, col V = "100 anos e mais") ---
$ws.Range("V1").Value = "100 anos e mais"

# --- New "Total" column header (W1) ---
$ws.Range("W1").Value = "Total"

# --- Row totals (column W) for existing category rows 2-6 ---
$ws.Range("W2").Value = 1907
$ws.Range("W3").Value = 186
$ws.Range("W4").Value = 990
$ws.Range("W5").Value = 274
$ws.Range("W6").Value = 1272

# --- New row 7: "Outros" category ---
$ws.Range("A7").Value = "Outros"

$row7 = New-Object 'object[,]' 1,21
$row7Values = @(169,7,13,48,81,95,103,80,102,97,105,102,100,100,156,127,162,142,86,28,6)
for ($i = 0; $i -lt 21; $i++) { $row7[0,$i] = $row7Values[$i] }
$ws.Range("B7:V7").Value = $row7

$ws.Range("W7").Value = 1909

# --- New row 8: "Total" row (column sums across rows 2-7) ---
$ws.Range("A8").Value = "Total"

$row8 = New-Object 'object[,]' 1,21
$row8Values = @(185,10,18,54,93,110,140,143,190,265,370,416,478,571,690,730,839,686,396,123,31)
for ($i = 0; $i -lt 21; $i++) { $row8[0,$i] = $row8Values[$i] }
$ws.Range("B8:V8").Value = $row8

$ws.Range("W8").Value = 6538
